$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "F-T1-G101-I401-CA",
    "F-T1-G101-I402-EP",
    "F-T1-G101-I403-CA",
    "F-T1-G101-I404-EP",
    "F-T1-G102-I405-CA",
    "F-T1-G102-I406-EP",
    "F-T1-G102-I407-CA",
    "F-T1-G102-I408-EP",
    "F-T2-G103-I409-CA",
    "F-T2-G103-I410-EP",
    "F-T2-G103-I411-CA",
    "F-T2-G103-I412-EP",
    "F-T2-G104-I413-CA",
    "F-T2-G104-I414-EP",
    "F-T2-G104-I415-CA",
    "F-T2-G104-I416-EP",
    "F-T3-G105-I417-CA",
    "F-T3-G105-I418-EP",
    "F-T3-G105-I419-CA",
    "F-T3-G105-I420-EP",
    "F-T3-G106-I421-CA",
    "F-T3-G106-I422-EP",
    "F-T3-G106-I423-CA",
    "F-T3-G106-I424-EP",
    "F-T4-G107-I425-CA",
    "F-T4-G107-I426-EP",
    "F-T4-G107-I427-CA",
    "F-T4-G107-I428-EP",
    "F-T4-G108-I429-CA",
    "F-T4-G108-I430-EP",
    "F-T4-G108-I431-CA",
    "F-T4-G108-I432-EP",
    "F-T5-G109-I433-CA",
    "F-T5-G109-I434-EP",
    "F-T5-G109-I435-CA",
    "F-T5-G109-I436-EP",
    "F-T5-G110-I437-CA",
    "F-T5-G110-I438-EP",
    "F-T5-G110-I439-CA",
    "F-T5-G110-I440-EP"
)

$startRow = 401
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
